$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) values for existing coin rows.
# Price values that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source data, which
# keeps formatted strings such as "27.808.66" or "0.0610" verbatim).
$ws.Range("D2").Value = "27.833.62"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.617.49"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'210.05"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'23.23"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "'0.0610"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "1.846.29"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "1.614.25"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "'0.558"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").Value = "'65.04"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "27.808.66"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'228.19"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "'10.05"
$ws.Range("E23").Value = "  -5.33%  "
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").Value = "'154.43"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "'15.44"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").Value = "'3.42"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").Value = "1.386.27"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").Value = "'0.992"
$ws.Range("E36").Value = "  +9.94%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("E40").Value = "  -4.58%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'0.995"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "'1.83"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "'65.33"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").Value = "1.757.74"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "'2.16"
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("E48").Value = "  -0.22%  "

# Row 49-51 shift: the BabyDogeCoin row is removed, Algorand and Cronos each
# move up one row, and a new coin (EnergySwap) is appended as the new row 51.
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.101"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0502"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.53"
$ws.Range("E51").Value = "  -1.04%  "

